$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, '3035 Campbell Place Aged Care Glen Waverley', 12),
    @(3, '3364 Assisi Centre Aged Care Rosanna', 21),
    @(4, '3622 Olivet Care Aged Care Services Ringwood', 13),
    @(5, '3633 Lifeview Emerald Glades Aged Care Emerald', 14),
    @(6, '3824 Estia Health South Morang', 11),
    @(7, '3961 Heritage Care Water Gardens Aged Care Facility Sydenham', 27),
    @(8, '4167 Royal Freemasons Centennial Lodge Wantirna South', 20),
    @(9, 'AG Industries Pty Ltd Factory Thomastown', 17),
    @(10, 'Aintree Primary School Aintree', 18),
    @(11, 'Australian Meat Group Abattoir Dandenong South', 16),
    @(12, 'Bacchus Marsh Childcare and Kindergarten Centre Bacchus Marsh', 30),
    @(13, 'Baden Powell College Tarneit', 15),
    @(14, 'Covenant College Bell Post Hill', 13),
    @(15, 'Dandenong South Primary School Dandenong', 11),
    @(16, 'Gladstone Park Secondary College 29 Oct Gladstone Park', 15),
    @(17, 'Hamlyn Banks Primary School Hamlyn Heights', 10),
    @(18, 'Hamlyn Views School Hamlyn Heights', 11),
    @(19, 'Hazelwood North Primary School Hazelwood North', 23),
    @(20, 'Hippity Hop Childcare and Kindergarten Pakenham', 10),
    @(21, 'M.C. Herd Corio', 10),
    @(22, 'Master Poultry Group West Footscray', 13),
    @(23, 'Morwell Park Primary School Morwell', 30),
    @(24, 'Nido Early School Woodend', 10),
    @(25, 'Northern Bay College Goldsworthy 9-12 Campus Corio', 18),
    @(26, 'Northern Bay College Wexford Campus Corio', 56),
    @(27, 'Northern Health Northern Hospital Epping Emergency Department Tier 1B', 32),
    @(28, 'Northern Health The Northern Hospital Epping', 10),
    @(29, 'Oakleigh South Primary School Oakleigh South', 16),
    @(30, 'Our Lady''s Catholic Primary School Wangaratta', 12),
    @(31, 'Sirius College Ibrahim Dellal Campus Sunshine', 11),
    @(32, 'Smartie Pants Early Learning and Development Diamond Creek', 19),
    @(33, 'St Georges Road Primary School Shepparton', 15),
    @(34, 'St Joseph''s School Quarry Hill', 27),
    @(35, 'St Josephs Catholic Primary School Warragul', 12),
    @(36, 'St Louis de Montfort''s School Aspendale', 12),
    @(37, 'St Vincents Hospital Emergency Department Melbourne', 14),
    @(38, 'Stockdale Road Primary School Traralgon', 31),
    @(39, 'Story House Early Learning Epping', 12),
    @(40, 'Sunbury Primary School Sunbury', 11),
    @(41, 'TUROSI PTY LTD Thomastown', 11),
    @(42, 'Templestowe Park Primary School Templestowe', 29),
    @(43, 'The Lake Primary School Cabarita', 14),
    @(44, 'Werribee Mercy Hospital Emergency Department', 25),
    @(45, 'Wodonga Cemetery Wodonga', 10),
    @(46, 'Wodonga Primary School Wodonga', 23),
    @(47, 'Wodonga Senior Secondary College Wodonga', 18),
    @(48, 'Wodonga South Primary School Wodonga', 28),
    @(49, 'Woodend Primary School Woodend', 20),
    @(50, 'Wyndham Christian College Wyndham Vale', 14)
)

foreach ($row in $data) {
    $r = $row[0]
    $name = $row[1]
    $val = $row[2]
    $ws.Cells.Item($r, 1).Value = $name
    $ws.Cells.Item($r, 2).Value = $val
}

